# Updated symbol list on Tue Dec 20 06:17:11 UTC 2022 with GitHub Actions
#
# The "Price" (D) and "Hora" (G) columns of the crypto price table are
# refreshed: most Price values get a small update and every Hora value
# moves from 5 to 6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep these columns stored as text (they already hold text-typed numeric
# strings like "249.20", "0.00000000750", "--", etc. in the workbook), so
# assigning plain numeric-looking strings doesn't get silently re-typed
# into floating point numbers and lose formatting / precision.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

# --- Price (column D) updates ---
$ws.Range("D2").Value = "248.96"
$ws.Range("D3").Value = "21.73"
$ws.Range("D4").Value = "5.364"
$ws.Range("D6").Value = "3.408"
$ws.Range("D7").Value = "6.386"
$ws.Range("D8").Value = "0.8142"
$ws.Range("D9").Value = "0.9609"
$ws.Range("D10").Value = "0.1417"
$ws.Range("D11").Value = "0.07595"
$ws.Range("D12").Value = "0.03195"
$ws.Range("D13").Value = "0.03049"
$ws.Range("D15").Value = "3.566"
$ws.Range("D16").Value = "0.001595"
$ws.Range("D17").Value = "0.04711"
$ws.Range("D18").Value = "0.0005766"
$ws.Range("D19").Value = "0.006237"
$ws.Range("D20").Value = "0.005086"
$ws.Range("D21").Value = "0.001033"
$ws.Range("D24").Value = "2.148"
$ws.Range("D25").Value = "0.3252"
$ws.Range("D28").Value = "0.0002998"
$ws.Range("D40").Value = "0.03948"
$ws.Range("D41").Value = "0.006981"
$ws.Range("D42").Value = "0.1067"
$ws.Range("D44").Value = "0.008769"
$ws.Range("D45").Value = "0.00005813"
$ws.Range("D47").Value = "0.0005496"
$ws.Range("D48").Value = "0.6795"
$ws.Range("D49").Value = "0.1694"

# --- Hora (column G) updates: every row 2..51 goes from "5" to "6" ---
$ws.Range("G2:G51").Value = "6"
